# "minor change to flow layout"
#
# Every top-level shape on the slide is nudged 200743 EMU to the left
# (vertical position, sizes, child offsets/extents inside groups, etc.
# all stay exactly as they were).
#
# PowerPoint's Shape.Left/.Top are expressed in points (single-precision)
# rather than EMU, and 200743 EMU is not a whole number of points
# (200743 / 12700 = 15.80654 pt). To make sure the value that ends up
# stored in the OOXML is exactly "old_emu - 200743" (and not off by a
# rounding unit because of the point<->EMU round trip), a tiny half-EMU
# epsilon is added before the Single-precision conversion truncates.

$EMU_PER_POINT = 12700.0
$deltaEmu = -200743
$deltaPt = $deltaEmu / $EMU_PER_POINT
$epsilonPt = 0.5 / $EMU_PER_POINT

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

foreach ($sh in $s.Shapes) {
    $sh.Left = $sh.Left + $deltaPt + $epsilonPt
}
